$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.445.53"
$ws.Range("D3").Value = "1.831.78"
$ws.Range("E3").Value = "  -1.97%  "
$ws.Range("E4").Value = "  -0.93%  "
$ws.Range("D5").Value = "'331.10"
$ws.Range("E5").Value = "  -0.71%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("D7").Value = "'0.4600"
$ws.Range("E7").Value = "  -2.65%  "
$ws.Range("D8").Value = "'0.3834"
$ws.Range("E8").Value = "  -1.81%  "
$ws.Range("D9").Value = "'46.58"
$ws.Range("E9").Value = "  -0.19%  "
$ws.Range("D10").Value = "'0.07887"
$ws.Range("E10").Value = "  -0.98%  "
$ws.Range("D11").Value = "'0.9720"
$ws.Range("E11").Value = "  -3.36%  "
$ws.Range("D12").Value = "'21.09"
$ws.Range("E12").Value = "  -2.20%  "
$ws.Range("D13").Value = "1.835.34"
$ws.Range("E13").Value = "  -2.20%  "
$ws.Range("D14").Value = "'5.885"
$ws.Range("E14").Value = "  -1.55%  "
$ws.Range("D15").Value = "'7.054"
$ws.Range("E15").Value = "  -1.21%  "
$ws.Range("D16").Value = "'1.003"
$ws.Range("E16").Value = "  -0.93%  "
$ws.Range("D17").Value = "'87.98"
$ws.Range("E17").Value = "  -0.18%  "
$ws.Range("D18").Value = "'0.06613"
$ws.Range("E18").Value = "  -1.31%  "
$ws.Range("D19").Value = "'0.00001030"
$ws.Range("E19").Value = "  -1.02%  "
$ws.Range("D20").Value = "'17.15"
$ws.Range("E20").Value = "  +1.29%  "
$ws.Range("E21").Value = "  -0.63%  "
$ws.Range("D22").Value = "27.456.22"
$ws.Range("E22").Value = "  -1.27%  "
$ws.Range("D23").Value = "'5.340"
$ws.Range("E23").Value = "  -2.21%  "
$ws.Range("E24").Value = "  -0.76%  "
$ws.Range("D25").Value = "'2.304"
$ws.Range("E25").Value = "  -1.05%  "
$ws.Range("D26").Value = "2.054.61"
$ws.Range("E26").Value = "  -2.08%  "
$ws.Range("D27").Value = "'157.08"
$ws.Range("E27").Value = "  -0.46%  "
$ws.Range("D28").Value = "'19.43"
$ws.Range("E28").Value = "  -1.32%  "
$ws.Range("D29").Value = "'2.062"
$ws.Range("E29").Value = "  -1.05%  "
$ws.Range("D30").Value = "'5.278"
$ws.Range("E30").Value = "  -2.35%  "
$ws.Range("D31").Value = "'118.77"
$ws.Range("E31").Value = "  -1.84%  "
$ws.Range("D32").Value = "'0.9553"
$ws.Range("E32").Value = "  -1.15%  "
$ws.Range("D33").Value = "'0.09288"
$ws.Range("E33").Value = "  -1.83%  "
$ws.Range("D34").Value = "'3.578"
$ws.Range("E34").Value = "  -1.43%  "
$ws.Range("D35").Value = "'5.235"
$ws.Range("E35").Value = "  -1.18%  "
$ws.Range("D36").Value = "'1.314"
$ws.Range("E36").Value = "  -2.05%  "
$ws.Range("D37").Value = "'0.05941"
$ws.Range("E37").Value = "  -1.50%  "
$ws.Range("D38").Value = "'0.02199"
$ws.Range("E38").Value = "  -0.80%  "
$ws.Range("D39").Value = "'8.046"
$ws.Range("E39").Value = "  -0.95%  "
$ws.Range("D40").Value = "'1.152"
$ws.Range("E40").Value = "  -4.26%  "
$ws.Range("D41").Value = "'0.5797"
$ws.Range("E41").Value = "  -1.97%  "
$ws.Range("D42").Value = "'0.1840"
$ws.Range("E42").Value = "  -2.41%  "
$ws.Range("D43").Value = "'10.00"
$ws.Range("E43").Value = "  -2.39%  "
$ws.Range("D44").Value = "'1.287"
$ws.Range("E44").Value = "  +2.54%  "
$ws.Range("D45").Value = "'0.5485"
$ws.Range("E45").Value = "  -2.37%  "
$ws.Range("D46").Value = "'11.94"
$ws.Range("E46").Value = "  -0.61%  "
$ws.Range("E47").Value = "  -2.03%  "
$ws.Range("D48").Value = "'0.06644"
$ws.Range("E48").Value = "  -1.95%  "
$ws.Range("D49").Value = "'110.41"
$ws.Range("E49").Value = "  -1.27%  "
$ws.Range("D50").Value = "'1.040"
$ws.Range("E50").Value = "  -2.13%  "
$ws.Range("D51").Value = "'1.001"
$ws.Range("E51").Value = "  -0.80%  "
